# "Generate Report for Handback" -- refresh the localization-status report:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - Latest Handback DateTime stamps are refreshed for zh-cn / de-de
#   - The stale "handback file is not latest" Error Detail warnings are cleared
#     now that the handback is in sync
#   - A couple of columns are re-sized (AutoFit-style) to reflect the new text

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) mirror the per-language
# sheets' Status column.
# ---------------------------------------------------------------------------
$overview = $wb.Sheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Widen the two status columns to fit the new, longer status text.
$overview.Range("E:E").ColumnWidth = 29.166666666666668
$overview.Range("F:F").ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-22 10:29:28"
$zhcn.Range("P2").Value = ""

$zhcn.Range("C:C").ColumnWidth = 29.166666666666668
$zhcn.Range("P:P").ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Sheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-22 10:29:35"
$dede.Range("P2").Value = ""

$dede.Range("C:C").ColumnWidth = 29.166666666666668
$dede.Range("P:P").ColumnWidth = 12.833333333333334
